$d = $word.ActiveDocument

# --- 1. Add the new "Abstract Title" paragraph style -----------------------
# wdStyleTypeParagraph = 1
$newStyle = $d.Styles.Add("Abstract Title", 1)

$newStyle.BaseStyle          = "Normal"
$newStyle.NextParagraphStyle = "Abstract"
$newStyle.QuickStyle         = $true

$newStyle.ParagraphFormat.KeepWithNext = $true
$newStyle.ParagraphFormat.KeepTogether = $true
# wdAlignParagraphCenter = 1
$newStyle.ParagraphFormat.Alignment    = 1
$newStyle.ParagraphFormat.SpaceAfter   = 0
$newStyle.ParagraphFormat.SpaceBefore  = 15

$newStyle.Font.Size   = 10
$newStyle.Font.SizeBi = 10
$newStyle.Font.Bold   = $true
$newStyle.Font.Color  = 9067060

# --- 2. Tighten the spacing above the existing "Abstract" style ------------
$abstract = $d.Styles.Item("Abstract")
$abstract.ParagraphFormat.SpaceBefore = 5
